$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '257.78'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.33%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '27.15'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-1.25%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.667'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-10.39%'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-0.59%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.642'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.52%'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-1.00%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9676'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-4.94%'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.75%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.03985'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '10.57%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07087'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-1.35%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03178'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-1.58%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09186'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.39%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001551'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.69%'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0006028'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.41%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006216'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '5.57%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.515'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.06%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.206'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-1.85%'
$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.204'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-0.90%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3079'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-2.27%'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-1.15%'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.856'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '9.05%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04228'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '1.24%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.08%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004304'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-4.72%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-0.04%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001936'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-0.14%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03828'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '0.09%'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1104'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-0.05%'
$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.003933'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-28.64%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002430'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-1.26%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01144'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '7.08%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005456'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '0.43%'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.08%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05998'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-45.05%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1300'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '5,877.54%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002099'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.08%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0001999'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.08%'
